$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1523.5536
$ws.Range("J17").Value = 1262.8334
$ws.Range("L17").Value = 3788.5002
$ws.Range("N17").Value = -4124.5002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 5698.5
$ws.Range("I28").Value = 265
$ws.Range("J28").Value = 21999
$ws.Range("K28").Value = 265
$ws.Range("L28").Value = 21999
$ws.Range("M28").Value = 220
$ws.Range("N28").Value = -22969

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1072.9429
$ws.Range("I40").Value = 1049.1333
$ws.Range("J40").Value = 1215.8
$ws.Range("K40").Value = 1049.1333
$ws.Range("L40").Value = 1215.8
$ws.Range("M40").Value = -874.1333
$ws.Range("N40").Value = -1565.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1145.1818
$ws.Range("J96").Value = 1433.3334
$ws.Range("L96").Value = 4300.0002
$ws.Range("N96").Value = -7046.0002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1499.325
$ws.Range("I112").Value = 900
$ws.Range("K112").Value = 2700
$ws.Range("M112").Value = -1592

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 554.5714
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 8773378
$ws.Range("I132").Value = 10418052
$ws.Range("K132").Value = 31254156
$ws.Range("M132").Value = -31251626

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3717.1384
$ws.Range("I32").Value = 2911.5344
$ws.Range("J32").Value = 10392.143
$ws.Range("K32").Value = 2911.5344
$ws.Range("L32").Value = 10392.143
$ws.Range("M32").Value = -2624.5344
$ws.Range("N32").Value = -10966.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1905.6
$ws.Range("I45").Value = 1932
$ws.Range("J45").Value = 1888
$ws.Range("K45").Value = 1932
$ws.Range("L45").Value = 1888
$ws.Range("M45").Value = -1555
$ws.Range("N45").Value = -2642

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 41668530
$ws.Range("I61").Value = 27778732
$ws.Range("J61").Value = 83337920
$ws.Range("K61").Value = 27778732
$ws.Range("L61").Value = 83337920
$ws.Range("M61").Value = -27778520
$ws.Range("N61").Value = -83338344

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1312.9048
$ws.Range("I74").Value = 952.6
$ws.Range("K74").Value = 952.6
$ws.Range("M74").Value = -78.60000000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1312.9048
$ws.Range("I77").Value = 952.6
$ws.Range("K77").Value = 4763
$ws.Range("M77").Value = -395

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 67886.5
$ws.Range("J109").Value = 67886.5
$ws.Range("L109").Value = 67886.5
$ws.Range("N109").Value = -70660.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1322.7587
$ws.Range("J110").Value = 2662.4
$ws.Range("L110").Value = 2662.4
$ws.Range("N110").Value = -6752.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1278.931
$ws.Range("I122").Value = 1221.3462
$ws.Range("K122").Value = 3664.0386
$ws.Range("M122").Value = -1214.0386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1377.9464
$ws.Range("I132").Value = 1027.2162
$ws.Range("J132").Value = 2060.9473
$ws.Range("K132").Value = 3081.6486
$ws.Range("L132").Value = 6182.841899999999
$ws.Range("M132").Value = -551.6486000000004
$ws.Range("N132").Value = -11242.8419

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 41668530
$ws.Range("I136").Value = 27778732
$ws.Range("J136").Value = 83337920
$ws.Range("K136").Value = 83336196
$ws.Range("L136").Value = 250013760
$ws.Range("M136").Value = -83333646
$ws.Range("N136").Value = -250018860

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 94029.95
$ws.Range("I86").Value = 3403.9333
$ws.Range("J86").Value = 288228.56
$ws.Range("K86").Value = 3403.9333
$ws.Range("L86").Value = 288228.56
$ws.Range("M86").Value = -2280.9333
$ws.Range("N86").Value = -290474.56

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 94029.95
$ws.Range("I89").Value = 3403.9333
$ws.Range("J89").Value = 288228.56
$ws.Range("K89").Value = 17019.6665
$ws.Range("L89").Value = 1441142.8
$ws.Range("M89").Value = -11403.6665
$ws.Range("N89").Value = -1452374.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1198.3125
$ws.Range("J107").Value = 1649.5
$ws.Range("L107").Value = 1649.5
$ws.Range("N107").Value = -5489.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4763.702
$ws.Range("I134").Value = 4547.5
$ws.Range("J134").Value = 6579.8
$ws.Range("K134").Value = 13642.5
$ws.Range("L134").Value = 19739.4
$ws.Range("M134").Value = -11107.5
$ws.Range("N134").Value = -24809.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2305865.5
$ws.Range("I31").Value = 5103280.5
$ws.Range("J31").Value = 2111.7646
$ws.Range("K31").Value = 5103280.5
$ws.Range("L31").Value = 2111.7646
$ws.Range("M31").Value = -5102985.5
$ws.Range("N31").Value = -2701.7646

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2305865.5
$ws.Range("I34").Value = 5103280.5
$ws.Range("J34").Value = 2111.7646
$ws.Range("K34").Value = 5103280.5
$ws.Range("L34").Value = 2111.7646
$ws.Range("M34").Value = -5103078.5
$ws.Range("N34").Value = -2515.7646

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3109533.5
$ws.Range("I58").Value = 8698612
$ws.Range("J58").Value = 4489.8887
$ws.Range("K58").Value = 8698612
$ws.Range("L58").Value = 4489.8887
$ws.Range("M58").Value = -8698409
$ws.Range("N58").Value = -4895.8887

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 50000
$ws.Range("J68").Value = 50000
$ws.Range("L68").Value = 50000
$ws.Range("N68").Value = -51498

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 50000
$ws.Range("J71").Value = 50000
$ws.Range("L71").Value = 150000
$ws.Range("N71").Value = -157488

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2060.25
$ws.Range("I132").Value = 1133.1904
$ws.Range("J132").Value = 4841.4287
$ws.Range("K132").Value = 3399.5712
$ws.Range("L132").Value = 14524.2861
$ws.Range("M132").Value = -869.5711999999999
$ws.Range("N132").Value = -19584.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1558.1632
$ws.Range("I134").Value = 1469.8108
$ws.Range("K134").Value = 4409.4324
$ws.Range("M134").Value = -1874.4324

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3109533.5
$ws.Range("I136").Value = 8698612
$ws.Range("J136").Value = 4489.8887
$ws.Range("K136").Value = 26095836
$ws.Range("L136").Value = 13469.6661
$ws.Range("M136").Value = -26093286
$ws.Range("N136").Value = -18569.6661

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 1870.3846
$ws.Range("I103").Value = 2381.75
$ws.Range("J103").Value = 1643.1111
$ws.Range("K103").Value = 7145.25
$ws.Range("L103").Value = 4929.3333
$ws.Range("M103").Value = -6266.25
$ws.Range("N103").Value = -6687.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 2936.125
$ws.Range("I114").Value = 700
$ws.Range("J114").Value = 3255.5715
$ws.Range("K114").Value = 2100
$ws.Range("L114").Value = 9766.7145
$ws.Range("M114").Value = 1154
$ws.Range("N114").Value = -16274.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 8594.134
$ws.Range("I131").Value = 404.2857
$ws.Range("J131").Value = 10102.789
$ws.Range("K131").Value = 1212.8571
$ws.Range("L131").Value = 30308.367
$ws.Range("M131").Value = 3827.1429
$ws.Range("N131").Value = -40388.367

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 360
$ws.Range("J107").Value = 500
$ws.Range("L107").Value = 500
$ws.Range("N107").Value = -4340

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 90000
$ws.Range("J110").Value = 90000
$ws.Range("L110").Value = 90000
$ws.Range("N110").Value = -98180

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2010.2
$ws.Range("I113").Value = 2350.3333
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 2350.3333
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -180.3332999999998
$ws.Range("N113").Value = -5840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1042024.4
$ws.Range("I132").Value = 1604285.8
$ws.Range("K132").Value = 4812857.4
$ws.Range("M132").Value = -4810327.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11122.8
$ws.Range("J40").Value = 9778.799999999999
$ws.Range("L40").Value = 9778.799999999999
$ws.Range("N40").Value = -10050.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2115.5557
$ws.Range("I68").Value = 1690
$ws.Range("J68").Value = 2966.6667
$ws.Range("K68").Value = 1690
$ws.Range("L68").Value = 2966.6667
$ws.Range("M68").Value = -941
$ws.Range("N68").Value = -4464.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2115.5557
$ws.Range("I71").Value = 1690
$ws.Range("J71").Value = 2966.6667
$ws.Range("K71").Value = 8450
$ws.Range("L71").Value = 14833.3335
$ws.Range("M71").Value = -4706
$ws.Range("N71").Value = -22321.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4935
$ws.Range("I122").Value = 6296.8335
$ws.Range("J122").Value = 3913.625
$ws.Range("K122").Value = 18890.5005
$ws.Range("L122").Value = 11740.875
$ws.Range("M122").Value = -16440.5005
$ws.Range("N122").Value = -16640.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1739.45
$ws.Range("I132").Value = 1122.4117
$ws.Range("K132").Value = 3367.2351
$ws.Range("M132").Value = -837.2351000000003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1896.0204
$ws.Range("I136").Value = 1155.7368
$ws.Range("K136").Value = 3467.2104
$ws.Range("M136").Value = -917.2103999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 29533.535
$ws.Range("I122").Value = 31703.424
$ws.Range("K122").Value = 95110.272
$ws.Range("M122").Value = -92660.272

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2918
$ws.Range("I126").Value = 1503.25
$ws.Range("K126").Value = 4509.75
$ws.Range("M126").Value = -2039.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1561.0714
$ws.Range("I132").Value = 1038.2142
$ws.Range("J132").Value = 3129.6428
$ws.Range("K132").Value = 3114.6426
$ws.Range("L132").Value = 9388.928400000001
$ws.Range("M132").Value = -584.6425999999997
$ws.Range("N132").Value = -14448.9284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 15874502
$ws.Range("I136").Value = 20577278
$ws.Range("J136").Value = 2636.125
$ws.Range("K136").Value = 61731834
$ws.Range("L136").Value = 7908.375
$ws.Range("M136").Value = -61729284
$ws.Range("N136").Value = -13008.375
